$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 148.35715
$ws.Range("I11").Value = 148.35715
$ws.Range("K11").Value = 148.35715
$ws.Range("M11").Value = -8.35714999999999
$ws.Range("H28").Value = 448.46667
$ws.Range("I28").Value = 409.35715
$ws.Range("K28").Value = 409.35715
$ws.Range("M28").Value = 75.64285000000001
$ws.Range("H40").Value = 5429.5
$ws.Range("I40").Value = 3921.5
$ws.Range("J40").Value = 6937.5
$ws.Range("K40").Value = 3921.5
$ws.Range("L40").Value = 6937.5
$ws.Range("M40").Value = -3746.5
$ws.Range("N40").Value = -7287.5
$ws.Range("H43").Value = 1353.2858
$ws.Range("I43").Value = 1393.25
$ws.Range("J43").Value = 1300
$ws.Range("K43").Value = 1393.25
$ws.Range("L43").Value = 1300
$ws.Range("M43").Value = -1324.25
$ws.Range("N43").Value = -1438
$ws.Range("H98").Value = 4021.7273
$ws.Range("I98").Value = 3999.558
$ws.Range("K98").Value = 3999.558
$ws.Range("M98").Value = -2501.558
$ws.Range("H107").Value = 454
$ws.Range("I107").Value = 392.77777
$ws.Range("K107").Value = 392.77777
$ws.Range("M107").Value = 1527.22223
$ws.Range("H111").Value = 2182.7778
$ws.Range("I111").Value = 2191
$ws.Range("K111").Value = 6573
$ws.Range("M111").Value = -3506
$ws.Range("H116").Value = 4085.7144
$ws.Range("H122").Value = 4021.7273
$ws.Range("I122").Value = 3999.558
$ws.Range("K122").Value = 11998.674
$ws.Range("M122").Value = -9548.673999999999
$ws.Range("H132").Value = 9094.846
$ws.Range("I132").Value = 9602.75
$ws.Range("K132").Value = 28808.25
$ws.Range("M132").Value = -26278.25
$ws.Range("H137").Value = 1936
$ws.Range("J137").Value = 2699
$ws.Range("L137").Value = 8097
$ws.Range("N137").Value = -13197

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5148.476
$ws.Range("J2").Value = 5122.7144
$ws.Range("L2").Value = 5122.7144
$ws.Range("N2").Value = -5348.7144
$ws.Range("H45").Value = 50282.254
$ws.Range("I45").Value = 69028.37
$ws.Range("K45").Value = 69028.37
$ws.Range("M45").Value = -68651.37
$ws.Range("H74").Value = 2635.818
$ws.Range("I74").Value = 1970.6
$ws.Range("K74").Value = 1970.6
$ws.Range("M74").Value = -1096.6
$ws.Range("H76").Value = 17678.75
$ws.Range("J76").Value = 17678.75
$ws.Range("L76").Value = 17678.75
$ws.Range("N76").Value = -18354.75
$ws.Range("H77").Value = 2635.818
$ws.Range("I77").Value = 1970.6
$ws.Range("K77").Value = 9853
$ws.Range("M77").Value = -5485
$ws.Range("H79").Value = 17678.75
$ws.Range("J79").Value = 17678.75
$ws.Range("L79").Value = 17678.75
$ws.Range("N79").Value = -20018.75
$ws.Range("H106").Value = 25370
$ws.Range("J106").Value = 25370
$ws.Range("L106").Value = 25370
$ws.Range("N106").Value = -27894
$ws.Range("H116").Value = 5148.476
$ws.Range("J116").Value = 5122.7144
$ws.Range("L116").Value = 5122.7144
$ws.Range("N116").Value = -9710.7144
$ws.Range("H132").Value = 4196.314
$ws.Range("I132").Value = 3352.2163
$ws.Range("K132").Value = 10056.6489
$ws.Range("M132").Value = -7526.6489

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5148.476
$ws.Range("J3").Value = 5122.7144
$ws.Range("L3").Value = 5122.7144
$ws.Range("N3").Value = -5350.7144
$ws.Range("H61").Value = 30000
$ws.Range("J61").Value = 30000
$ws.Range("L61").Value = 30000
$ws.Range("N61").Value = -30626
$ws.Range("H64").Value = 492.81818
$ws.Range("J64").Value = 266
$ws.Range("L64").Value = 266
$ws.Range("N64").Value = -716
$ws.Range("H67").Value = 492.81818
$ws.Range("J67").Value = 266
$ws.Range("L67").Value = 266
$ws.Range("N67").Value = -1826
$ws.Range("H99").Value = 4802.4
$ws.Range("J99").Value = 6442.2
$ws.Range("L99").Value = 6442.2
$ws.Range("N99").Value = -9438.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4158
$ws.Range("I31").Value = 3132.1155
$ws.Range("J31").Value = 6380.75
$ws.Range("K31").Value = 3132.1155
$ws.Range("L31").Value = 6380.75
$ws.Range("M31").Value = -2837.1155
$ws.Range("N31").Value = -6970.75
$ws.Range("H34").Value = 4158
$ws.Range("I34").Value = 3132.1155
$ws.Range("J34").Value = 6380.75
$ws.Range("K34").Value = 3132.1155
$ws.Range("L34").Value = 6380.75
$ws.Range("M34").Value = -2930.1155
$ws.Range("N34").Value = -6784.75
$ws.Range("H74").Value = 42281.668
$ws.Range("J74").Value = 42281.668
$ws.Range("L74").Value = 42281.668
$ws.Range("N74").Value = -44029.668
$ws.Range("H77").Value = 42281.668
$ws.Range("J77").Value = 42281.668
$ws.Range("L77").Value = 126845.004
$ws.Range("N77").Value = -135581.004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1791.2727
$ws.Range("I5").Value = 529.1429
$ws.Range("J5").Value = 4000
$ws.Range("K5").Value = 1587.4287
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = -1475.4287
$ws.Range("N5").Value = -12224
$ws.Range("H8").Value = 680.2
$ws.Range("I8").Value = 680.2
$ws.Range("K8").Value = 2040.6
$ws.Range("M8").Value = -1901.6
$ws.Range("H131").Value = 4547345
$ws.Range("J131").Value = 7695093.5
$ws.Range("L131").Value = 23085280.5
$ws.Range("N131").Value = -23095360.5
$ws.Range("H134").Value = 1702.5454
$ws.Range("I134").Value = 1702.5454
$ws.Range("K134").Value = 5107.6362
$ws.Range("M134").Value = -37.63619999999992
$ws.Range("H135").Value = 1791.2727
$ws.Range("I135").Value = 529.1429
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 4762.2861
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -2227.2861
$ws.Range("N135").Value = -41070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H80").Value = 3783.6
$ws.Range("I80").Value = 3471.25
$ws.Range("J80").Value = 5033
$ws.Range("K80").Value = 3471.25
$ws.Range("L80").Value = 5033
$ws.Range("M80").Value = -2473.25
$ws.Range("N80").Value = -7029
$ws.Range("H83").Value = 3783.6
$ws.Range("I83").Value = 3471.25
$ws.Range("J83").Value = 5033
$ws.Range("K83").Value = 17356.25
$ws.Range("L83").Value = 25165
$ws.Range("M83").Value = -12364.25
$ws.Range("N83").Value = -35149
$ws.Range("H97").Value = 9858.333
$ws.Range("I97").Value = 1883.3334
$ws.Range("K97").Value = 1883.3334
$ws.Range("M97").Value = -1387.3334
$ws.Range("H113").Value = 2000750
$ws.Range("I113").Value = 2000750
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000750
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -1998580
$ws.Range("H132").Value = 7401.174
$ws.Range("I132").Value = 7211.4
$ws.Range("J132").Value = 8666.333
$ws.Range("K132").Value = 21634.2
$ws.Range("L132").Value = 25998.999
$ws.Range("M132").Value = -19104.2
$ws.Range("N132").Value = -31058.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2262.7742
$ws.Range("I16").Value = 2193.913
$ws.Range("K16").Value = 2193.913
$ws.Range("M16").Value = -2023.913
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H63").Value = 45973.75
$ws.Range("J63").Value = 45945
$ws.Range("L63").Value = 45945
$ws.Range("N63").Value = -47443
$ws.Range("H66").Value = 45973.75
$ws.Range("J66").Value = 45945
$ws.Range("L66").Value = 137835
$ws.Range("N66").Value = -145323
$ws.Range("H93").Value = 24211.555
$ws.Range("I93").Value = 1129.2858
$ws.Range("K93").Value = 1129.2858
$ws.Range("M93").Value = 118.7141999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 4277.3335
$ws.Range("I6").Value = 3916
$ws.Range("K6").Value = 3916
$ws.Range("M6").Value = -3801
$ws.Range("H100").Value = 864.8571
$ws.Range("J100").Value = 420
$ws.Range("L100").Value = 840
$ws.Range("N100").Value = -1922
$ws.Range("H113").Value = 284.5625
$ws.Range("I113").Value = 300.53333
$ws.Range("K113").Value = 901.5999899999999
$ws.Range("M113").Value = 1268.40001
$ws.Range("H122").Value = 3797.2942
$ws.Range("I122").Value = 2864.7778
$ws.Range("J122").Value = 4846.375
$ws.Range("K122").Value = 8594.3334
$ws.Range("L122").Value = 14539.125
$ws.Range("M122").Value = -6144.3334
$ws.Range("N122").Value = -19439.125
